$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 2296.88
$ws.Range("J6").Value = 465.91666
$ws.Range("L6").Value = 1397.74998
$ws.Range("N6").Value = -1621.74998

# Row 21
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -532

# Row 23
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -766

# Row 34
$ws.Range("H34").Value = 16000
$ws.Range("I34").Value = 16000
$ws.Range("K34").Value = 16000
$ws.Range("M34").Value = -15797

# Row 36
$ws.Range("H36").Value = 16000
$ws.Range("I36").Value = 16000
$ws.Range("K36").Value = 16000
$ws.Range("M36").Value = -15285

# Row 47
$ws.Range("H47").Value = 2000
$ws.Range("I47").Value = 1000
$ws.Range("J47").Value = 3000
$ws.Range("K47").Value = 1000
$ws.Range("L47").Value = 3000
$ws.Range("M47").Value = -28
$ws.Range("N47").Value = -4944

# Row 74
$ws.Range("H74").Value = 3632.5334
$ws.Range("I74").Value = 3235.2727
$ws.Range("J74").Value = 4725
$ws.Range("K74").Value = 3235.2727
$ws.Range("L74").Value = 4725
$ws.Range("M74").Value = -2299.2727
$ws.Range("N74").Value = -6597

# Row 77
$ws.Range("H77").Value = 3632.5334
$ws.Range("I77").Value = 3235.2727
$ws.Range("J77").Value = 4725
$ws.Range("K77").Value = 16176.3635
$ws.Range("L77").Value = 23625
$ws.Range("M77").Value = -11496.3635
$ws.Range("N77").Value = -32985

# Row 82
$ws.Range("H82").Value = 293.66666
$ws.Range("I82").Value = 293.66666
$ws.Range("K82").Value = 880.9999799999999
$ws.Range("M82").Value = -474.9999799999999

# Row 85
$ws.Range("H85").Value = 293.66666
$ws.Range("I85").Value = 293.66666
$ws.Range("K85").Value = 880.9999799999999
$ws.Range("M85").Value = 523.0000200000001

# Row 88
$ws.Range("H88").Value = 2505.7144
$ws.Range("I88").Value = 2003
$ws.Range("J88").Value = 2589.5
$ws.Range("K88").Value = 2003
$ws.Range("L88").Value = 2589.5
$ws.Range("M88").Value = -1597
$ws.Range("N88").Value = -3401.5

# Row 91
$ws.Range("H91").Value = 2505.7144
$ws.Range("I91").Value = 2003
$ws.Range("J91").Value = 2589.5
$ws.Range("K91").Value = 2003
$ws.Range("L91").Value = 2589.5
$ws.Range("M91").Value = -599
$ws.Range("N91").Value = -5397.5

# Row 107
$ws.Range("H107").Value = 8541.049999999999
$ws.Range("I107").Value = 10902.667
$ws.Range("K107").Value = 10902.667
$ws.Range("M107").Value = -8982.666999999999

# Row 132
$ws.Range("H132").Value = 4640.276
$ws.Range("I132").Value = 4021.923
$ws.Range("K132").Value = 12065.769
$ws.Range("M132").Value = -9535.769

# Row 140
$ws.Range("H140").Value = 80694.8
$ws.Range("J140").Value = 80694.8
$ws.Range("L140").Value = 80694.8
$ws.Range("N140").Value = -91054.8

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 550
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 900
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = -84
$ws.Range("N4").Value = -1132

# Row 32
$ws.Range("H32").Value = 3201.9744
$ws.Range("I32").Value = 3266.946
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 3266.946
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -2979.946
$ws.Range("N32").Value = -2574

# Row 74
$ws.Range("H74").Value = 2579
$ws.Range("I74").Value = 1391.1351
$ws.Range("K74").Value = 1391.1351
$ws.Range("M74").Value = -517.1351

# Row 77
$ws.Range("H77").Value = 2579
$ws.Range("I77").Value = 1391.1351
$ws.Range("K77").Value = 6955.675499999999
$ws.Range("M77").Value = -2587.675499999999

# Row 139
$ws.Range("H139").Value = 74499.5
$ws.Range("J139").Value = 74499.5
$ws.Range("L139").Value = 74499.5
$ws.Range("N139").Value = -84779.5

$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10336

# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# Row 26
$ws.Range("H26").Value = 29996.5
$ws.Range("I26").Value = 29996.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 29996.5
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -29704.5

# Row 86
$ws.Range("H86").Value = 4802.778
$ws.Range("J86").Value = 3000.6924
$ws.Range("L86").Value = 3000.6924
$ws.Range("N86").Value = -5246.6924

# Row 89
$ws.Range("H89").Value = 4802.778
$ws.Range("J89").Value = 3000.6924
$ws.Range("L89").Value = 15003.462
$ws.Range("N89").Value = -26235.462

# Row 134
$ws.Range("H134").Value = 2829.2144
$ws.Range("I134").Value = 1966.091
$ws.Range("K134").Value = 5898.272999999999
$ws.Range("M134").Value = -3363.272999999999

# Row 141
$ws.Range("H141").Value = 73666.664
$ws.Range("J141").Value = 73666.664
$ws.Range("L141").Value = 73666.664
$ws.Range("N141").Value = -84026.664

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 6002
$ws.Range("I3").Value = 6002
$ws.Range("K3").Value = 6002
$ws.Range("M3").Value = -5889

# Row 16
$ws.Range("H16").Value = 1238.6428
$ws.Range("I16").Value = 1149.0834
$ws.Range("K16").Value = 1149.0834
$ws.Range("M16").Value = -862.0834

# Row 22
$ws.Range("H22").Value = 813.5714
$ws.Range("I22").Value = 565
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 565
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -215
$ws.Range("N22").Value = -1700

# Row 69
$ws.Range("H69").Value = 10097
$ws.Range("I69").Value = 10097
$ws.Range("K69").Value = 10097
$ws.Range("M69").Value = -9348

# Row 72
$ws.Range("H72").Value = 10097
$ws.Range("I72").Value = 10097
$ws.Range("K72").Value = 30291
$ws.Range("M72").Value = -26547

# Row 113
$ws.Range("H113").Value = 1238.6428
$ws.Range("I113").Value = 1149.0834
$ws.Range("K113").Value = 1149.0834
$ws.Range("M113").Value = 1020.9166

# Row 134
$ws.Range("H134").Value = 5673.0415
$ws.Range("I134").Value = 4811.524
$ws.Range("K134").Value = 14434.572
$ws.Range("M134").Value = -11899.572

# Row 140
$ws.Range("H140").Value = 94997.5
$ws.Range("J140").Value = 94997.5
$ws.Range("L140").Value = 94997.5
$ws.Range("N140").Value = -105357.5

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1564540.5
$ws.Range("I34").Value = 2779194.5
$ws.Range("J34").Value = 2842.7144
$ws.Range("K34").Value = 8337583.5
$ws.Range("L34").Value = 8528.143199999999
$ws.Range("M34").Value = -8337499.5
$ws.Range("N34").Value = -8696.143199999999

# Row 56
$ws.Range("H56").Value = 6605.5713
$ws.Range("I56").Value = 6605.5713
$ws.Range("K56").Value = 6605.5713
$ws.Range("M56").Value = -6075.5713

# Row 109
$ws.Range("H109").Value = 2742.6667
$ws.Range("I109").Value = 1827.2858
$ws.Range("J109").Value = 5946.5
$ws.Range("K109").Value = 5481.857400000001
$ws.Range("L109").Value = 17839.5
$ws.Range("M109").Value = -4441.857400000001
$ws.Range("N109").Value = -19919.5

# Row 113
$ws.Range("H113").Value = 1768.875
$ws.Range("J113").Value = 1858.1666
$ws.Range("L113").Value = 5574.4998
$ws.Range("N113").Value = -9914.4998

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 81.521736
$ws.Range("I2").Value = 81.521736
$ws.Range("K2").Value = 81.521736
$ws.Range("M2").Value = 31.478264

# Row 19
$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5576

# Row 33
$ws.Range("H33").Value = 19989.5
$ws.Range("J33").Value = 19989.5
$ws.Range("L33").Value = 19989.5
$ws.Range("N33").Value = -20493.5

# Row 48
$ws.Range("H48").Value = 18142.857
$ws.Range("I48").Value = 10000
$ws.Range("J48").Value = 19500
$ws.Range("K48").Value = 10000
$ws.Range("L48").Value = 19500
$ws.Range("M48").Value = -9515
$ws.Range("N48").Value = -20470

# Row 52
$ws.Range("H52").Value = 31000
$ws.Range("J52").Value = 31000
$ws.Range("L52").Value = 31000
$ws.Range("N52").Value = -31518

$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 10015000

# Row 93
$ws.Range("H93").Value = 4563.926
$ws.Range("I93").Value = 4508.6924
$ws.Range("K93").Value = 4508.6924
$ws.Range("M93").Value = -3260.6924

# Row 136
$ws.Range("H136").Value = 10530.625
$ws.Range("I136").Value = 39500
$ws.Range("K136").Value = 118500
$ws.Range("M136").Value = -115950

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 11199.375
$ws.Range("I2").Value = 11199.375
$ws.Range("K2").Value = 11199.375
$ws.Range("M2").Value = -11087.375

# Row 107
$ws.Range("H107").Value = 32574.35
$ws.Range("I107").Value = 3106.2856
$ws.Range("K107").Value = 9318.856800000001
$ws.Range("M107").Value = -7398.856800000001

# Row 132
$ws.Range("H132").Value = 81205.5
$ws.Range("I132").Value = 81205.5
$ws.Range("K132").Value = 243616.5
$ws.Range("M132").Value = -241086.5

# Row 136
$ws.Range("H136").Value = 1865.9375
$ws.Range("I136").Value = 1281.4783
$ws.Range("J136").Value = 3359.5557
$ws.Range("K136").Value = 3844.4349
$ws.Range("L136").Value = 10078.6671
$ws.Range("M136").Value = -1294.4349
$ws.Range("N136").Value = -15178.6671
